# repull data, push all data, mean calculation
# Update column F ("dSF") values to reflect newly pulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -3
    4  = 6
    5  = -1
    7  = -3
    9  = 2
    10 = 1
    11 = 2
    13 = 1
    14 = -3
    15 = -3
    16 = -2
    17 = 1
    20 = 3
    22 = 7
    23 = -1
    24 = -5
    25 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
